$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current PI Fringe / insurance subsidy figure bumped.
$ws.Range("I10").Value = 1965

# Benefits/Fringe for PI#2 (row 8) now reference the "Current PI Fringe
# Amount" cell ($I$24) instead of the old hard-coded 16.64% constant.
$ws.Range("D8").Formula = '=D7*$I$24'
$ws.Range("E8").Formula = '=E7*$I$24'

# Reflect the author's new cursor position / selection on the sheet.
$ws.Range("I12").Select()

$wb.Application.Calculate()
